$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Flag column (F) for rows 3-10 from "Yes" to "No"
for ($r = 3; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = "No"
}

# Update the view: top-left cell and active-cell selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("H7").Select()
